$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.460.18'
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.150.09'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '613.16'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.71'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.146.91'
$ws.Range('D8').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.35'
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.30%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.46'
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.674.99'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.118'
$ws.Range('D16').Style = "Normal"

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.451.73'
$ws.Range('D17').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.153.15'
$ws.Range('D18').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.84'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '476.41'
$ws.Range('D20').Style = "Normal"

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.65'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.71'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.27'
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.37%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.59'
$ws.Range('D28').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.43'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +8.76%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.09'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -5.10%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.52'
$ws.Range('D33').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.93'
$ws.Range('D36').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.62%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0748'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.12%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.09'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +4.66%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '455.38'
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0396'
$ws.Range('D41').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.32'
$ws.Range('D43').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.854.86'
$ws.Range('D44').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.43'
$ws.Range('D47').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +4.89%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '26.41'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Arweave'
$ws.Range('B51').Style = "Normal"

$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('C51').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '34.04'
$ws.Range('D51').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +4.03%  '
$ws.Range('E51').Style = "Normal"
